# Auto commit update: refresh Metrics source figures and the "today"
# sheet's comparison formulas, then leave the workbook's active
# sheet/selection the way the author left it (Metrics!C8 selected,
# "today" scrolled to F11:F22).

$wb = $excel.ActiveWorkbook
$wsMetrics = $wb.Worksheets.Item("Metrics")
$wsToday   = $wb.Worksheets.Item("today")

# --- Metrics!B2:B13 — updated source metrics -----------------------------
$wsMetrics.Range("B2").Value  = 210212.69
$wsMetrics.Range("B3").Value  = 157289.44
$wsMetrics.Range("B4").Value  = 51069.320000000007
$wsMetrics.Range("B5").Value  = 8524
$wsMetrics.Range("B6").Value  = 5846083.419999999
$wsMetrics.Range("B7").Value  = 4928007.07
$wsMetrics.Range("B8").Value  = 1715161.14
$wsMetrics.Range("B9").Value  = 228801
$wsMetrics.Range("B10").Value = 34311464.409999996
$wsMetrics.Range("B11").Value = 32203282.23
$wsMetrics.Range("B12").Value = 11996883.18
$wsMetrics.Range("B13").Value = 1326431

# --- today!B3:B6 — newly filled-in daily adjustment formulas -------------
$wsToday.Range("B3").Formula = "=13186.83+4627.69"
$wsToday.Range("B4").Formula = "=9019.52+3892.22"
$wsToday.Range("B5").Formula = "=2918.73+1313.1"
$wsToday.Range("B6").Formula = "=541+176"

# --- view state: "today" keeps its own selection, Metrics becomes the ---
# --- active sheet with C8 selected (matches the captured sheetViews) ----
$wsToday.Range("F11:F22").Select()
$wsMetrics.Activate()
$wsMetrics.Range("C8").Select()
